$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1018.3158
$ws.Range("I2").Value = 251.76923
$ws.Range("J2").Value = 2679.1667
$ws.Range("K2").Value = 251.76923
$ws.Range("L2").Value = 2679.1667
$ws.Range("M2").Value = -138.76923
$ws.Range("N2").Value = -2905.1667
$ws.Range("H9").Value = 468.15
$ws.Range("I9").Value = 143.92857
$ws.Range("J9").Value = 1224.6666
$ws.Range("K9").Value = 143.92857
$ws.Range("L9").Value = 1224.6666
$ws.Range("M9").Value = 25.07142999999999
$ws.Range("N9").Value = -1562.6666
$ws.Range("H33").Value = 951.7273
$ws.Range("I33").Value = 1108.8889
$ws.Range("K33").Value = 1108.8889
$ws.Range("M33").Value = -879.8888999999999
$ws.Range("H38").Value = 763
$ws.Range("I38").Value = 292.22223
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 876.66669
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -504.66669
$ws.Range("N38").Value = -15744
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -431
$ws.Range("N43").ClearContents()
$ws.Range("H58").Value = 1247.6666
$ws.Range("I58").Value = 121.5
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 364.5
$ws.Range("L58").Value = 10500
$ws.Range("M58").Value = -214.5
$ws.Range("N58").Value = -10800
$ws.Range("H62").Value = 5195.875
$ws.Range("I62").Value = 2858
$ws.Range("J62").Value = 6598.6
$ws.Range("K62").Value = 2858
$ws.Range("L62").Value = 6598.6
$ws.Range("M62").Value = -2234
$ws.Range("N62").Value = -7846.6
$ws.Range("H64").Value = 5333.3335
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H65").Value = 5195.875
$ws.Range("I65").Value = 2858
$ws.Range("J65").Value = 6598.6
$ws.Range("K65").Value = 14290
$ws.Range("L65").Value = 32993
$ws.Range("M65").Value = -11170
$ws.Range("N65").Value = -39233
$ws.Range("H67").Value = 5333.3335
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H100").Value = 2900.8
$ws.Range("I100").Value = 2376
$ws.Range("K100").Value = 2376
$ws.Range("M100").Value = -1835
$ws.Range("H137").Value = 2045.7
$ws.Range("I137").Value = 1917.8334
$ws.Range("K137").Value = 5753.5002
$ws.Range("M137").Value = -3203.5002
$ws.Range("H138").Value = 1933.2325
$ws.Range("J138").Value = 2908
$ws.Range("L138").Value = 8724
$ws.Range("N138").Value = -19004
$ws.Range("H141").Value = 1394.375
$ws.Range("I141").Value = 1425.2174
$ws.Range("J141").Value = 685
$ws.Range("K141").Value = 4275.6522
$ws.Range("L141").Value = 2055
$ws.Range("M141").Value = 904.3477999999996
$ws.Range("N141").Value = -12415

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1583.96
$ws.Range("I74").Value = 1545.7916
$ws.Range("K74").Value = 1545.7916
$ws.Range("M74").Value = -671.7916
$ws.Range("H77").Value = 1583.96
$ws.Range("I77").Value = 1545.7916
$ws.Range("K77").Value = 7728.958000000001
$ws.Range("M77").Value = -3360.958000000001
$ws.Range("H122").Value = 1804.7778
$ws.Range("I122").Value = 799
$ws.Range("J122").Value = 3816.3333
$ws.Range("K122").Value = 2397
$ws.Range("L122").Value = 11448.9999
$ws.Range("M122").Value = 53
$ws.Range("N122").Value = -16348.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2478.8572
$ws.Range("I134").Value = 2478.8572
$ws.Range("K134").Value = 7436.571599999999
$ws.Range("M134").Value = -4901.571599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 86.611115
$ws.Range("J7").Value = 141
$ws.Range("L7").Value = 141
$ws.Range("N7").Value = -367
$ws.Range("H122").Value = 1692.3125
$ws.Range("I122").Value = 1692.3125
$ws.Range("K122").Value = 5076.9375
$ws.Range("M122").Value = -2626.9375
$ws.Range("H134").Value = 937.2857
$ws.Range("I134").Value = 937.2857
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2811.8571
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -276.8571000000002
$ws.Range("N134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 387.72223
$ws.Range("I92").Value = 341.57144
$ws.Range("J92").Value = 549.25
$ws.Range("K92").Value = 1024.71432
$ws.Range("L92").Value = 1647.75
$ws.Range("M92").Value = 223.28568
$ws.Range("N92").Value = -4143.75
$ws.Range("H98").Value = 2371.8572
$ws.Range("J98").Value = 4001.5
$ws.Range("L98").Value = 12004.5
$ws.Range("N98").Value = -15000.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 953.7857
$ws.Range("I102").Value = 953.7857
$ws.Range("K102").Value = 953.7857
$ws.Range("M102").Value = 668.2143
$ws.Range("H122").Value = 1869.1
$ws.Range("J122").Value = 4100
$ws.Range("L122").Value = 12300
$ws.Range("N122").Value = -17200
$ws.Range("H126").Value = 2078
$ws.Range("I126").Value = 1829
$ws.Range("J126").Value = 2949.5
$ws.Range("K126").Value = 5487
$ws.Range("L126").Value = 8848.5
$ws.Range("M126").Value = -3017
$ws.Range("N126").Value = -13788.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6947.095
$ws.Range("I7").Value = 3700
$ws.Range("J7").Value = 8570.643
$ws.Range("K7").Value = 3700
$ws.Range("L7").Value = 8570.643
$ws.Range("M7").Value = -3588
$ws.Range("N7").Value = -8794.643
$ws.Range("H22").Value = 1999.909
$ws.Range("I22").Value = 1899.9
$ws.Range("K22").Value = 1899.9
$ws.Range("M22").Value = -1604.9
$ws.Range("H27").Value = 1999.909
$ws.Range("I27").Value = 1899.9
$ws.Range("K27").Value = 1899.9
$ws.Range("M27").Value = -1792.9
$ws.Range("H100").Value = 2633.4546
$ws.Range("I100").Value = 2675.5715
$ws.Range("K100").Value = 2675.5715
$ws.Range("M100").Value = -2134.5715
$ws.Range("H126").Value = 6947.095
$ws.Range("I126").Value = 3700
$ws.Range("J126").Value = 8570.643
$ws.Range("K126").Value = 11100
$ws.Range("L126").Value = 25711.929
$ws.Range("M126").Value = -8630
$ws.Range("N126").Value = -30651.929
$ws.Range("H132").Value = 2357.1538
$ws.Range("I132").Value = 2189.2
$ws.Range("K132").Value = 6567.599999999999
$ws.Range("M132").Value = -4037.599999999999
$ws.Range("H139").Value = 64049
$ws.Range("I139").Value = 64049
$ws.Range("K139").Value = 64049
$ws.Range("M139").Value = -58909

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 8334248.5
$ws.Range("I100").Value = 8334248.5
$ws.Range("K100").Value = 16668497
$ws.Range("M100").Value = -16667956
$ws.Range("H132").Value = 2248.2
$ws.Range("I132").Value = 2250.2273
$ws.Range("J132").Value = 2233.3333
$ws.Range("K132").Value = 6750.6819
$ws.Range("L132").Value = 6699.999899999999
$ws.Range("M132").Value = -4220.6819
$ws.Range("N132").Value = -11759.9999
